# Lac_AllBounds13.xlsx - "Added Flow vs R1L to the cell data modeled by tissue slice code"
#
# Mirrors the existing Kpl-vs-R1L summary block (rows 22:24, built off column B)
# with a new Flow_Lac-vs-R1L summary block (rows 38:40, built off column F):
#   row 38 -> group headers (HK-2 / UMRC6 / UOK262 / UOK + DIDS) in G:J
#   row 39 -> row label "Flow_Lac" in F, AVERAGE() of each group's F-column data in G:J
#   row 40 -> standard error (STDEV/SQRT(COUNT)) of each group's F-column data in G:J

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Best-effort: give the workbook/sheet their VBA-style code names (the engine may not
# persist this, but it costs nothing to ask for it).
try { $wb.CodeName = "ThisWorkbook" } catch {}
try { $ws.CodeName = "Sheet1" } catch {}

# --- Row 38: group headers, reusing the same labels/shared strings as B22:F22 ---
$ws.Range("G38").Value = "HK-2"
$ws.Range("H38").Value = "UMRC6"
$ws.Range("I38").Value = "UOK262"
$ws.Range("J38").Value = "UOK + DIDS"

# --- Row 39: series label + group averages of column F ---
$ws.Range("F39").Value = "Flow_Lac"
$ws.Range("G39").Formula = "=AVERAGE(F`$1:F`$3)"
$ws.Range("H39").Formula = "=AVERAGE(F`$4:F`$6)"
$ws.Range("I39").Formula = "=AVERAGE(F`$9:F`$11)"
$ws.Range("J39").Formula = "=AVERAGE(F`$13:F`$16)"

# --- Row 40: group standard errors of column F ---
$ws.Range("G40").Formula = "=STDEV(F`$1:F`$3)/SQRT(COUNT(F`$1:F`$3))"
$ws.Range("H40").Formula = "=STDEV(F`$4:F`$6)/SQRT(COUNT(F`$4:F`$6))"
$ws.Range("I40").Formula = "=STDEV(F`$9:F`$11)/SQRT(COUNT(F`$9:F`$11))"
$ws.Range("J40").Formula = "=STDEV(F`$13:F`$16)/SQRT(COUNT(F`$13:F`$16))"

# Scroll the view down near the new block and select it, matching the saved
# workbook view (topLeftCell="A16", selection F38:J40).
try { $excel.ActiveWindow.ScrollRow = 16 } catch {}
$ws.Range("F38:J40").Select()
